# Apply updated cryptocurrency price/volume data to columns D (Price) and E (Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.700.64'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '3.586.33'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''609.26'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '''145.84'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''0.493'
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("D9").Value = '''0.137'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '''7.98'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").Value = '4.185.24'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").Value = '''0.0000209'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '''30.15'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = '3.588.18'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '66.746.33'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '''11.44'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = '''15.11'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").Value = '''433.38'
$ws.Range("E21").Value = '  +1.08%  '
$ws.Range("D22").Value = '''0.621'
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("D23").Value = '''79.33'
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("D24").Value = '3.723.94'
$ws.Range("E24").Value = '  +0.83%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  -2.13%  '
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").Value = '''2.52'
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("D29").Value = '''9.23'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").Value = '3.574.42'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("D33").Value = '''25.48'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '''0.157'
$ws.Range("E34").Value = '  -1.79%  '
$ws.Range("D35").Value = '''7.88'
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("D38").Value = '''5.65'
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("D39").Value = '''173.41'
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").Value = '''0.0855'
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("D41").Value = '''5.24'
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("D42").Value = '''0.893'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").Value = '''1.95'
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '''2.54'
$ws.Range("E45").Value = '  +5.54%  '
$ws.Range("E46").Value = '  -1.90%  '
$ws.Range("D47").Value = '''25.22'
$ws.Range("E47").Value = '  -3.24%  '
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("D49").Value = '''23.69'
$ws.Range("E49").Value = '  +2.84%  '
$ws.Range("D50").Value = '''0.945'
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").Value = '''0.237'
$ws.Range("E51").Value = '  -1.01%  '
